$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" ----
$ovw = $wb.Worksheets.Item("Overview")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"

# Latest HO Xliff Generate Date
$ovw.Range("G2").Value = "2016-10-26 08:32:42"
$ovw.Range("G3").Value = "2016-10-26 08:32:42"

# Column widths for zh-cn (E) and de-de (F) columns
$ovw.Columns.Item(5).ColumnWidth = 16.3333333333333
$ovw.Columns.Item(6).ColumnWidth = 16.3333333333333

# ---- Sheet "zh-cn" ----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"

$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

$zhcn.Range("H2").Value = "2016-10-26 08:32:30"
$zhcn.Range("H3").Value = "2016-10-26 08:32:30"

$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cfe77ee5a4dc650768a299f54012b62f2f25504b/e2e/b3830289-c780-410f-9b55-a9e2659232bc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b71da43590173647fae294ae9a85b89ee203a592/e2e/b3830289-c780-410f-9b55-a9e2659232bc.md."

$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- Sheet "de-de" ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"

$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

$dede.Range("H2").Value = "2016-10-26 08:32:42"
$dede.Range("H3").Value = "2016-10-26 08:32:42"

$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cfe77ee5a4dc650768a299f54012b62f2f25504b/e2e/b3830289-c780-410f-9b55-a9e2659232bc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b71da43590173647fae294ae9a85b89ee203a592/e2e/b3830289-c780-410f-9b55-a9e2659232bc.md."

$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
